$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("questions") - update question topics for Q3 and Q7
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B8").Value = "you could choose the next Binance listing"
$ws1.Range("B4").Value = "one of them was more likely to launch a rug"

$ws1.Activate()
$ws1.Range("B7").Select()

# ---------------------------------------------------------------------------
# Sheet 2 ("answers") - add a new pair of total-votes options for Q3
# (Ronaldo / Sophie Rain / Livvy Dune / Mr. Beast / Ishowspeed), replacing
# the old "superpower" options (Invisibility / Flying / Time Travel).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# The old Q3 block lives in rows 28:30. Keep the first 3 rows in place and
# insert 2 more rows right after them so the Q3 block grows from 3 to 5 rows.
$ws2.Rows("31:32").Insert()

# Row 28 - Ronaldo
$ws2.Range("A28").Value = "Q3"
$ws2.Range("B28").Value = "Ronaldo"
$ws2.Range("C28").Value = "cristiano"
$ws2.Range("D28").Value = "jpg"

# Row 29 - Sophie Rain
$ws2.Range("A29").Value = "Q3"
$ws2.Range("B29").Value = "Sophie Rain"
$ws2.Range("C29").Value = "sophie"
$ws2.Range("D29").Value = "jpg"

# Row 30 - Livvy Dune
$ws2.Range("A30").Value = "Q3"
$ws2.Range("B30").Value = "Livvy Dune"
$ws2.Range("C30").Value = "livvy"
$ws2.Range("D30").Value = "jpg"

# Row 31 - Mr. Beast (new row)
$ws2.Range("A31").Value = "Q3"
$ws2.Range("B31").Value = "Mr. Beast"
$ws2.Range("C31").Value = "mrBeast"
$ws2.Range("D31").Value = "jpg"

# Row 32 - Ishowspeed (new row)
$ws2.Range("A32").Value = "Q3"
$ws2.Range("C32").Value = "speed"
$ws2.Range("B32").Value = "Ishowspeed"
$ws2.Range("D32").Value = "jpg"

# Rebuild every formula in column E (image path lookup) so that:
#  - the stray literal value that had crept into row 27 becomes a real
#    formula again, and
#  - the "/images/" prefix used everywhere becomes "images/" (no leading
#    slash), matching the site's updated asset paths.
$lastRow = $ws2.Range("A" + $ws2.Rows.Count).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws2.Range("E" + $r).Formula = '=_xlfn.CONCAT("images/",C' + $r + ',".",D' + $r + ')'
}

$ws2.Activate()
$ws2.Range("C31").Select()
